$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-29 Sunday", "2024-09-30 Monday"),
    @("53÷6=", "71÷9="),
    @("78÷7=", "16÷6="),
    @("63÷4=", "98÷6="),
    @("77÷9=", "65÷8="),
    @("25÷6=", "13÷5="),
    @("54÷4=", "41÷8="),
    @("32÷6=", "64÷5="),
    @("49÷8=", "93÷5="),
    @("76÷2=", "95÷6="),
    @("47÷3=", "59÷5="),
    @("20÷3=", "16÷5="),
    @("13÷8=", "74÷5="),
    @("26÷8=", "50÷8="),
    @("63÷8=", "65÷6="),
    @("86÷9=", "63÷5="),
    @("35÷4=", "63÷6="),
    @("28÷5=", "15÷5="),
    @("77÷5=", "44÷3="),
    @("11÷3=", "63÷9="),
    @("10÷5=", "31÷8="),
    @("28÷3=", "89÷8="),
    @("66÷7=", "63÷7="),
    @("88÷3=", "19÷4="),
    @("69÷2=", "47÷5="),
    @("99÷5=", "94÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
